# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 07:32"

# Rows 27 and 28 swap identity (Israel now ranks above Ucrania) and get refreshed data.
# Row 27 becomes Israel with newly updated figures.
$ws.Range("A27").Value = "Israel"
$ws.Range("B27").Value = 146542
$ws.Range("C27").Value = 1016
$ws.Range("D27").Value = 111539
$ws.Range("E27").Value = 33926
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 1077

# Row 28 becomes Ucrania, carrying its previous (unchanged) figures.
$ws.Range("A28").Value = "Ucrania"
$ws.Range("B28").Value = 145612
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 65877
$ws.Range("E28").Value = 76712
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 3023

# Row 61 (Uzbekistan) refreshed figures
$ws.Range("B61").Value = 45679
$ws.Range("C61").Value = 206
$ws.Range("D61").Value = 42555
$ws.Range("E61").Value = 2751
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 373

# Row 160 (Belice) refreshed figures
$ws.Range("B160").Value = 1399
$ws.Range("C160").Value = 34
$ws.Range("D160").Value = 360
$ws.Range("E160").Value = 1020
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 19

# Row 172 (Islas Turcas y Caicos) refreshed figures
$ws.Range("B172").Value = 638
$ws.Range("C172").Value = 10
$ws.Range("D172").Value = 270
$ws.Range("E172").Value = 363

# Row 209 (Laos) refreshed figures
$ws.Range("B209").Value = 23
$ws.Range("C209").Value = 1
$ws.Range("E209").Value = 2
